$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percent/price cells are stored as plain text in this workbook (values like
# "3.105.77" are thousands-grouped strings, not numbers). Force Text format
# before writing so Excel does not auto-coerce numeric-looking strings (e.g.
# "580.16") into floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.067.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.105.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.61%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.96"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.78%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.100.93"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.51%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.43"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.37%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.481"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.617.87"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.071.67"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.105.86"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.26"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.54"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.56"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.39"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.97%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.03"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.18%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.08"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.42"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.68"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.98"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.67%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.92"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.991"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.73"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.12"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.98%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.31%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.77%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.846.48"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0362"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "385.01"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.07"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.76%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.06"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.84%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.89%  "
